$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.379.32"
$ws.Range("E2").Value = "  -2.96%  "

# Row 3
$ws.Range("D3").Value = "1.748.67"
$ws.Range("E3").Value = "  -3.73%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'321.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.42%  "

# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

# Row 7
$ws.Range("D7").Value = "'0.4228"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.71%  "

# Row 8
$ws.Range("D8").Value = "'0.3597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.72%  "

# Row 9
$ws.Range("D9").Value = "'0.07514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.30%  "

# Row 10
$ws.Range("D10").Value = "'42.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.89%  "

# Row 11
$ws.Range("D11").Value = "'1.089"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.22%  "

# Row 12
$ws.Range("D12").Value = "'1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13
$ws.Range("D13").Value = "'20.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.64%  "

# Row 14
$ws.Range("D14").Value = "'6.024"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.86%  "

# Row 15
$ws.Range("D15").Value = "'7.198"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.92%  "

# Row 16
$ws.Range("D16").Value = "1.754.08"
$ws.Range("E16").Value = "  -5.22%  "

# Row 17
$ws.Range("D17").Value = "'91.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.23%  "

# Row 18
$ws.Range("D18").Value = "'0.00001066"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.71%  "

# Row 19
$ws.Range("D19").Value = "'0.06349"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.77%  "

# Row 20
$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "

# Row 21
$ws.Range("D21").Value = "'17.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.93%  "

# Row 22
$ws.Range("D22").Value = "'5.874"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.04%  "

# Row 23
$ws.Range("D23").Value = "27.419.80"
$ws.Range("E23").Value = "  -3.06%  "

# Row 24
$ws.Range("D24").Value = "'11.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.30%  "

# Row 25
$ws.Range("D25").Value = "'2.085"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.55%  "

# Row 26
$ws.Range("D26").Value = "'160.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.00%  "

# Row 27
$ws.Range("D27").Value = "'20.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.54%  "

# Row 28
$ws.Range("D28").Value = "1.955.29"
$ws.Range("E28").Value = "  -3.89%  "

# Row 29
$ws.Range("D29").Value = "'2.123"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.51%  "

# Row 30
$ws.Range("D30").Value = "'123.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.81%  "

# Row 31
$ws.Range("D31").Value = "'1.108"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.26%  "

# Row 32
$ws.Range("D32").Value = "'3.644"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.44%  "

# Row 33
$ws.Range("D33").Value = "'5.539"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.58%  "

# Row 34
$ws.Range("D34").Value = "'0.08825"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.14%  "

# Row 35
$ws.Range("D35").Value = "'12.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.46%  "

# Row 36
$ws.Range("D36").Value = "'0.02275"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.29%  "

# Row 37
$ws.Range("D37").Value = "'0.2094"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.91%  "

# Row 38
$ws.Range("E38").Value = "  -3.45%  "

# Row 39
$ws.Range("D39").Value = "'0.6310"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.92%  "

# Row 40
$ws.Range("D40").Value = "'4.933"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.22%  "

# Row 41
$ws.Range("D41").Value = "'1.175"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.80%  "

# Row 42
$ws.Range("D42").Value = "'1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "

# Row 43
$ws.Range("D43").Value = "'7.843"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.90%  "

# Row 44
$ws.Range("D44").Value = "'1.389"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "

# Row 45
$ws.Range("D45").Value = "'13.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.49%  "

# Row 46
$ws.Range("D46").Value = "'0.5852"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.63%  "

# Row 47
$ws.Range("D47").Value = "'3.687"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.967"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.49%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'122.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.11%  "

# Row 50
$ws.Range("D50").Value = "'1.150"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "

# Row 51
$ws.Range("D51").Value = "'0.06794"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
